$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are not
# auto-converted to numbers by Excel type inference on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.410.14"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "3.686.73"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "418.89"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "129.92"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").Value = "3.679.45"
$ws.Range("E7").Value = "  +4.48%  "
$ws.Range("D8").Value = "0.641"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "0.766"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").Value = "0.179"
$ws.Range("E11").Value = "  +7.79%  "
$ws.Range("D12").Value = "0.0000390"
$ws.Range("E12").Value = "  +44.52%  "
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "10.60"
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "4.275.62"
$ws.Range("E15").Value = "  +4.52%  "
$ws.Range("D16").Value = "0.139"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "20.54"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "3.677.85"
$ws.Range("E18").Value = "  +4.09%  "
$ws.Range("D19").Value = "13.28"
$ws.Range("E19").Value = "  +5.41%  "
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "66.488.53"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "443.17"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "16.42"
$ws.Range("E23").Value = "  +22.49%  "
$ws.Range("D24").Value = "89.91"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "37.25"
$ws.Range("E26").Value = "  +8.09%  "
$ws.Range("D27").Value = "10.19"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "5.13"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("E30").Value = "  +8.59%  "
$ws.Range("D31").Value = "12.69"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "2.72"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").Value = "7.31"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").Value = "41.63"
$ws.Range("D36").Value = "57.28"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  +33.91%  "
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "29.33"
$ws.Range("E42").Value = "  +33.33%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").Value = "148.97"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "4.37"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.88"
$ws.Range("E49").Value = "  -7.51%  "
$ws.Range("D50").Value = "0.307"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("E51").Value = "  +13.11%  "

# Restore the default (General) style on column D so no stray number format
# is left applied to the cells (matches original inline-string formatting).
$ws.Range("D2:D51").Style = "Normal"
